$wb = $excel.ActiveWorkbook

# --- Spring 1: leave data untouched, just move the cursor off the cell it
# was sitting on before the author switched over to Spring 2 ---
$ws1 = $wb.Worksheets.Item("Spring 1")
$ws1.Activate()
$ws1.Range("C8").Select()

# --- Spring 2: fill in the "Avance" column (C) for HU3, then frame the
# whole data block with a thin box border ---
$ws2 = $wb.Worksheets.Item("Spring 2")
$ws2.Activate()

$ws2.Range("C2").Value = 6
$ws2.Range("C3").Value = 5
$ws2.Range("C4").Value = 4
$ws2.Range("C5").Value = 3
$ws2.Range("C6").Value = 2
$ws2.Range("C7").Value = 1

$tableRange = $ws2.Range("A2:C7")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# Leave the cursor where the author ended up after data entry / formatting
$ws2.Range("C14").Select()
